$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1145185.2
$ws.Range("I9").Value = 259.66666
$ws.Range("K9").Value = 259.66666
$ws.Range("M9").Value = -90.66665999999998
$ws.Range("H12").Value = 566.6667
$ws.Range("I12").Value = 350
$ws.Range("K12").Value = 350
$ws.Range("M12").Value = -180
$ws.Range("H32").Value = 8000
$ws.Range("J32").Value = 7999.5
$ws.Range("L32").Value = 7999.5
$ws.Range("N32").Value = -8651.5
$ws.Range("H38").Value = 970.4286
$ws.Range("J38").Value = 2166.5
$ws.Range("L38").Value = 6499.5
$ws.Range("N38").Value = -7243.5
$ws.Range("H40").Value = 3752.9285
$ws.Range("I40").Value = 3552.1428
$ws.Range("J40").Value = 4355.2856
$ws.Range("K40").Value = 3552.1428
$ws.Range("L40").Value = 4355.2856
$ws.Range("M40").Value = -3377.1428
$ws.Range("N40").Value = -4705.2856
$ws.Range("H74").Value = 7600
$ws.Range("I74").Value = 5885.7144
$ws.Range("K74").Value = 5885.7144
$ws.Range("M74").Value = -4949.7144
$ws.Range("H77").Value = 7600
$ws.Range("I77").Value = 5885.7144
$ws.Range("K77").Value = 29428.572
$ws.Range("M77").Value = -24748.572
$ws.Range("H100").Value = 6239.9844
$ws.Range("I100").Value = 1547.6522
$ws.Range("J100").Value = 8809.596
$ws.Range("K100").Value = 1547.6522
$ws.Range("L100").Value = 8809.596
$ws.Range("M100").Value = -1006.6522
$ws.Range("N100").Value = -9891.596
$ws.Range("H103").Value = 1140.2354
$ws.Range("I103").Value = 1088.5834
$ws.Range("J103").Value = 1264.2
$ws.Range("K103").Value = 3265.7502
$ws.Range("L103").Value = 3792.6
$ws.Range("M103").Value = -2679.7502
$ws.Range("N103").Value = -4964.6
$ws.Range("H138").Value = 2805.3635
$ws.Range("J138").Value = 3846.0667
$ws.Range("L138").Value = 11538.2001
$ws.Range("N138").Value = -21818.2001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 10071.866
$ws.Range("I45").Value = 10700.692
$ws.Range("K45").Value = 10700.692
$ws.Range("M45").Value = -10323.692
$ws.Range("H74").Value = 199473.8
$ws.Range("I74").Value = 215450.3
$ws.Range("J74").Value = 11750
$ws.Range("K74").Value = 215450.3
$ws.Range("L74").Value = 11750
$ws.Range("M74").Value = -214576.3
$ws.Range("N74").Value = -13498
$ws.Range("H77").Value = 199473.8
$ws.Range("I77").Value = 215450.3
$ws.Range("J77").Value = 11750
$ws.Range("K77").Value = 1077251.5
$ws.Range("L77").Value = 58750
$ws.Range("M77").Value = -1072883.5
$ws.Range("N77").Value = -67486
$ws.Range("H102").Value = 1395.5555
$ws.Range("I102").Value = 1395.5555
$ws.Range("K102").Value = 1395.5555
$ws.Range("M102").Value = 226.4445000000001
$ws.Range("H141").Value = 119984
$ws.Range("J141").Value = 119984
$ws.Range("L141").Value = 119984
$ws.Range("N141").Value = -130344

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 6250.25
$ws.Range("I22").Value = 4000.5
$ws.Range("K22").Value = 4000.5
$ws.Range("M22").Value = -3827.5
$ws.Range("H82").Value = 27330.4
$ws.Range("I82").Value = 16200.5
$ws.Range("J82").Value = 127499.5
$ws.Range("K82").Value = 16200.5
$ws.Range("L82").Value = 127499.5
$ws.Range("M82").Value = -15817.5
$ws.Range("N82").Value = -128265.5
$ws.Range("H85").Value = 27330.4
$ws.Range("I85").Value = 16200.5
$ws.Range("J85").Value = 127499.5
$ws.Range("K85").Value = 16200.5
$ws.Range("L85").Value = 127499.5
$ws.Range("M85").Value = -14874.5
$ws.Range("N85").Value = -130151.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2692.9412
$ws.Range("I31").Value = 1718.2894
$ws.Range("J31").Value = 5541.923
$ws.Range("K31").Value = 1718.2894
$ws.Range("L31").Value = 5541.923
$ws.Range("M31").Value = -1423.2894
$ws.Range("N31").Value = -6131.923
$ws.Range("H34").Value = 2692.9412
$ws.Range("I34").Value = 1718.2894
$ws.Range("J34").Value = 5541.923
$ws.Range("K34").Value = 1718.2894
$ws.Range("L34").Value = 5541.923
$ws.Range("M34").Value = -1516.2894
$ws.Range("N34").Value = -5945.923
$ws.Range("H58").Value = 6206
$ws.Range("I58").Value = 4404
$ws.Range("J58").Value = 8308.333000000001
$ws.Range("K58").Value = 4404
$ws.Range("L58").Value = 8308.333000000001
$ws.Range("M58").Value = -4201
$ws.Range("N58").Value = -8714.333000000001
$ws.Range("H122").Value = 9735.409
$ws.Range("I122").Value = 11608.588
$ws.Range("K122").Value = 34825.764
$ws.Range("M122").Value = -32375.764
$ws.Range("H136").Value = 6206
$ws.Range("I136").Value = 4404
$ws.Range("J136").Value = 8308.333000000001
$ws.Range("K136").Value = 13212
$ws.Range("L136").Value = 24924.999
$ws.Range("M136").Value = -10662
$ws.Range("N136").Value = -30024.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1061.3846
$ws.Range("I5").Value = 725.5
$ws.Range("J5").Value = 1598.8
$ws.Range("K5").Value = 2176.5
$ws.Range("L5").Value = 4796.4
$ws.Range("M5").Value = -2064.5
$ws.Range("N5").Value = -5020.4
$ws.Range("H35").Value = 1857.1428
$ws.Range("H107").Value = 7960.125
$ws.Range("J107").Value = 7960.125
$ws.Range("L107").Value = 23880.375
$ws.Range("N107").Value = -27720.375
$ws.Range("H135").Value = 1061.3846
$ws.Range("I135").Value = 725.5
$ws.Range("J135").Value = 1598.8
$ws.Range("K135").Value = 6529.5
$ws.Range("L135").Value = 14389.2
$ws.Range("M135").Value = -3994.5
$ws.Range("N135").Value = -19459.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 21918.613
$ws.Range("I122").Value = 19809.764
$ws.Range("K122").Value = 59429.292
$ws.Range("M122").Value = -56979.292
$ws.Range("H123").Value = 24999
$ws.Range("J123").Value = 24999
$ws.Range("L123").Value = 24999
$ws.Range("N123").Value = -29899

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2609.7693
$ws.Range("I16").Value = 2620.5
$ws.Range("J16").Value = 2592.6
$ws.Range("K16").Value = 2620.5
$ws.Range("L16").Value = 2592.6
$ws.Range("M16").Value = -2450.5
$ws.Range("N16").Value = -2932.6
$ws.Range("H22").Value = 3521.6765
$ws.Range("I22").Value = 2011.421
$ws.Range("J22").Value = 5434.6665
$ws.Range("K22").Value = 2011.421
$ws.Range("L22").Value = 5434.6665
$ws.Range("M22").Value = -1716.421
$ws.Range("N22").Value = -6024.6665
$ws.Range("H27").Value = 3521.6765
$ws.Range("I27").Value = 2011.421
$ws.Range("J27").Value = 5434.6665
$ws.Range("K27").Value = 2011.421
$ws.Range("L27").Value = 5434.6665
$ws.Range("M27").Value = -1904.421
$ws.Range("N27").Value = -5648.6665
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("H68").Value = 2789.05
$ws.Range("J68").Value = 3350.875
$ws.Range("L68").Value = 3350.875
$ws.Range("N68").Value = -4848.875
$ws.Range("H71").Value = 2789.05
$ws.Range("J71").Value = 3350.875
$ws.Range("L71").Value = 16754.375
$ws.Range("N71").Value = -24242.375
$ws.Range("H133").Value = 67330.664
$ws.Range("J133").Value = 67330.664
$ws.Range("L133").Value = 67330.664
$ws.Range("N133").Value = -72390.664
$ws.Range("H141").Value = 71712.5
$ws.Range("J141").Value = 71712.5
$ws.Range("L141").Value = 71712.5
$ws.Range("N141").Value = -82072.5
$ws.Range("N43").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 74999
$ws.Range("J75").Value = 74999
$ws.Range("L75").Value = 74999
$ws.Range("N75").Value = -76871
$ws.Range("H78").Value = 74999
$ws.Range("J78").Value = 74999
$ws.Range("L78").Value = 224997
$ws.Range("N78").Value = -234357
$ws.Range("H132").Value = 5586
$ws.Range("I132").Value = 4726.6924
$ws.Range("J132").Value = 13033.333
$ws.Range("K132").Value = 14180.0772
$ws.Range("L132").Value = 39099.999
$ws.Range("M132").Value = -11650.0772
$ws.Range("N132").Value = -44159.999
$ws.Range("H136").Value = 1995.5667
$ws.Range("I136").Value = 1013.5926
$ws.Range("K136").Value = 3040.7778
$ws.Range("M136").Value = -490.7777999999998
